$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "Abbotsford"
$ws.Cells.Item(2,2).Value = "Bodriggy Brewing Company  245 Johnston Street, Abbotsford VIC 3067"
$ws.Cells.Item(2,3).Value = "28/12/20 2:50pm-5:30pm"
$ws.Cells.Item(2,4).Value = "Case dined at venue"

$ws.Cells.Item(3,1).Value = "Albert Park"
$ws.Cells.Item(3,2).Value = "The Guilty Moose Cafe  143 Victoria Avenue, Albert Park VIC 3206"
$ws.Cells.Item(3,3).Value = "30/12/20 1pm-1:30pm"
$ws.Cells.Item(3,4).Value = "Case ate at cafe"

$ws.Cells.Item(4,1).Value = "Bairnsdale"
$ws.Cells.Item(4,2).Value = "V/Line train - Bairnsdale to Melbourne"
$ws.Cells.Item(4,3).Value = "30/12/20 12:45pm-4:30pm"
$ws.Cells.Item(4,4).Value = "Case caught train from Bairnsdale to Caulfield"

$ws.Cells.Item(5,1).Value = "Brighton"
$ws.Cells.Item(5,2).Value = "Sons of Mary Restaurant  14 Spink St, Brighton VIC 3186"
$ws.Cells.Item(5,3).Value = "24/12/20 10:00am-11:05am"
$ws.Cells.Item(5,4).Value = "Case ate at restaurant"

$ws.Cells.Item(6,1).Value = "Camberwell"
$ws.Cells.Item(6,2).Value = "Crown Nails  766 Riversdale Road, Camberwell VIC 3124"
$ws.Cells.Item(6,3).Value = "30/12/20 1:30pm-2:30pm"
$ws.Cells.Item(6,4).Value = "Case attended"

$ws.Cells.Item(7,1).Value = "Camberwell"
$ws.Cells.Item(7,2).Value = "Tao Dumplings  1 Evans Place, Camberwell VIC 3124"
$ws.Cells.Item(7,3).Value = "29/12/20 12:30pm-1:30pm"
$ws.Cells.Item(7,4).Value = "Case ate at restaurant"

$ws.Cells.Item(8,1).Value = "Caulfield"
$ws.Cells.Item(8,2).Value = "Metro Train - Frankston line"
$ws.Cells.Item(8,3).Value = "30/12/20 4:30pm-5:00pm"
$ws.Cells.Item(8,4).Value = "Case caught train from Caulfield to Cheltenham"

$ws.Cells.Item(9,1).Value = "Cheltenham"
$ws.Cells.Item(9,2).Value = "Two Bob Snob, 256 Charman Road"
$ws.Cells.Item(9,3).Value = "22/12/2020 1pm - 2pm"
$ws.Cells.Item(9,4).Value = "Case attended Venue"

$ws.Cells.Item(10,1).Value = "Docklands"
$ws.Cells.Item(10,2).Value = "Melbourne Boat Hire - Yarra River Cruise. 45 Newquay Promenade, Docklands VIC 3008"
$ws.Cells.Item(10,3).Value = "28/12/2020 11.26am-2:00pm"
$ws.Cells.Item(10,4).Value = "Case attended venue"

$ws.Cells.Item(11,1).Value = "Doveton"
$ws.Cells.Item(11,2).Value = "Holy Family Parish Doveton Catholic  100 Power Road, Doveton VIC 3177"
$ws.Cells.Item(11,3).Value = "26/12/20 6:30pm"
$ws.Cells.Item(11,4).Value = "Case attended English service"

$ws.Cells.Item(12,1).Value = "Hampton"
$ws.Cells.Item(12,2).Value = "Merrymen Cafe, 2 Small Street, Hampton VIC"
$ws.Cells.Item(12,3).Value = "28-12-2020 1:00pm-2:00pm"
$ws.Cells.Item(12,4).Value = "Case ate in store"

$ws.Cells.Item(13,1).Value = "Hampton"
$ws.Cells.Item(13,2).Value = "Merrymen Cafe, 2 Small Street, Hampton VIC"
$ws.Cells.Item(13,3).Value = "28-12-2020 1:30pm-2:30pm"
$ws.Cells.Item(13,4).Value = "Case ate in store"

$ws.Cells.Item(14,1).Value = "Lakes Entrance"
$ws.Cells.Item(14,2).Value = "Albert and Co. Cafe - Bellevue Hotel  201 Esplanade, Lakes Entrance VIC 3909"
$ws.Cells.Item(14,3).Value = "29/12/20 09:30am-10:45am"
$ws.Cells.Item(14,4).Value = "Case ate in cafe"

$ws.Cells.Item(15,1).Value = "Lakes Entrance"
$ws.Cells.Item(15,2).Value = "Chants Summer Carnival - Footbridge, Lakes Entrance VIC 3909"
$ws.Cells.Item(15,3).Value = "29/12/20 7:00pm-9:30pm"
$ws.Cells.Item(15,4).Value = "Case attended carnival"

$ws.Cells.Item(16,1).Value = "Lakes Entrance"
$ws.Cells.Item(16,2).Value = "Esplanade Resort Lakes Entrance - Hotel bar  1 Esplanade, Lakes Entrance VIC 3909"
$ws.Cells.Item(16,3).Value = "29/12/20 05:30pm-05:50pm"
$ws.Cells.Item(16,4).Value = "Case in hotel bar"

$ws.Cells.Item(17,1).Value = "Lakes Entrance"
$ws.Cells.Item(17,2).Value = "Esplanade Resort Lakes Entrance - Pool area  1 Esplanade, Lakes Entrance VIC 3909"
$ws.Cells.Item(17,3).Value = "29/12/20 03:30pm-04:30pm"
$ws.Cells.Item(17,4).Value = "Case visited pool area"

$ws.Cells.Item(18,1).Value = "Lakes Entrance"
$ws.Cells.Item(18,2).Value = "Lakes Boat Shed Cafe  54 Marine Parade, Lakes Entrance VIC 3909"
$ws.Cells.Item(18,3).Value = "30/12/20 9:30am-10:30am"
$ws.Cells.Item(18,4).Value = "Case attended cafe"

$ws.Cells.Item(19,1).Value = "Lakes Entrance"
$ws.Cells.Item(19,2).Value = "The Esplanade Resort and Spa  1 Esplanade, Lakes Entrance VIC 3909"
$ws.Cells.Item(19,3).Value = "29/12/20 2:30pm-5:50pm"
$ws.Cells.Item(19,4).Value = "Case attended site"

$ws.Cells.Item(20,1).Value = "Lakes Entrance"
$ws.Cells.Item(20,2).Value = "V/Line bus - Lakes Entrance to Bairnsdale"
$ws.Cells.Item(20,3).Value = "30/12/20 11:55am-12:30pm"
$ws.Cells.Item(20,4).Value = "Case caught the 11:55am bus from Lakes Entrance"

$ws.Cells.Item(21,1).Value = "Lakes Entrance"
$ws.Cells.Item(21,2).Value = "Wyanga Winery  248 Baades Rd, Lakes Entrance VIC 3909"
$ws.Cells.Item(21,3).Value = "29/12/20 1:00pm-2:00pm"
$ws.Cells.Item(21,4).Value = "Case visited vineyard"

$ws.Cells.Item(22,1).Value = "Leongatha"
$ws.Cells.Item(22,2).Value = "Coral Fish and Chips 53 Bair St, Leongatha VIC 3953"
$ws.Cells.Item(22,3).Value = "26/12/20 5:30pm-7:30pm"
$ws.Cells.Item(22,4).Value = "Case worked in store"

$ws.Cells.Item(23,1).Value = "Leongatha"
$ws.Cells.Item(23,2).Value = "Coral Fish and Chips 53 Bair St, Leongatha VIC 3953"
$ws.Cells.Item(23,3).Value = "27/12/20 4:00pm-7:30pm"
$ws.Cells.Item(23,4).Value = "Case worked in store"

$ws.Cells.Item(24,1).Value = "Melbourne"
$ws.Cells.Item(24,2).Value = "European Bier Cafe City  120 Exhibition Street Melbourne VIC 3000"
$ws.Cells.Item(24,3).Value = "28/12/20 8:00pm-9:30pm"
$ws.Cells.Item(24,4).Value = "Case attended cafe"

$ws.Cells.Item(25,1).Value = "Melbourne"
$ws.Cells.Item(25,2).Value = "Fonda Mexican Flinders Lane  31 Flinders Lane Melbourne"
$ws.Cells.Item(25,3).Value = "29/12/20 6:00pm-7:30pm"
$ws.Cells.Item(25,4).Value = "Case attended restaurant"

$ws.Cells.Item(26,1).Value = "Melbourne"
$ws.Cells.Item(26,2).Value = "Melbourne Central Lion Hotel, 211 La Trobe Street"
$ws.Cells.Item(26,3).Value = "28/12/2020 10pm - 12.30am"
$ws.Cells.Item(26,4).Value = "Case attended venue"

$ws.Cells.Item(27,1).Value = "Moorabin"
$ws.Cells.Item(27,2).Value = "Grape and Grain Liquor Cellars, 14/16 Station St"
$ws.Cells.Item(27,3).Value = "21/12/20 2pm - 10pm  22/12/20 10am - 6pm  24/12/20 1pm - 10pm  28/12/20 8.05pm - 8.47pm  29/12/20 12pm - 4pm"
$ws.Cells.Item(27,4).Value = "Cases workplace"

$ws.Cells.Item(28,1).Value = "Oakleigh"
$ws.Cells.Item(28,2).Value = "Melissa Oakleigh Restaurant  6 Eaton Mall, Oakleigh VIC 3166"
$ws.Cells.Item(28,3).Value = "28/12/20 7:30pm-8:15pm"
$ws.Cells.Item(28,4).Value = "Case dined in restaurant"

$ws.Cells.Item(29,1).Value = "Sandringham Line"
$ws.Cells.Item(29,2).Value = "Metro Train line Sandringham"
$ws.Cells.Item(29,3).Value = "28/12/20 7pm -7.50pm"
$ws.Cells.Item(29,4).Value = "Travelled by train from Sandringham Station to Parliament Station"
